# REACH.xlsx brochure config update: refresh the promotions table rows
# (row 2: tuition fees brochure, row 3: offshore Q4 promotions brochure)
# with the new file names / download links, and leave the selection on
# the last edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells(2, 1).Value = "REACH Tuition Course Fees 2021.pdf"
$ws.Cells(3, 1).Value = "REACH Offshore Q4 Promotions Region 1.pdf"

$ws.Cells(2, 2).Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/seapae/reach/Reach_Tuition_Course_Fees_2021_v1.0.pdf"
$ws.Cells(3, 2).Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/seapae/reach/Reach-Offshore-SEAPAE-Q4-Promotions-1OCT-31DEC21_VOL-1.1.pdf"

[void]$ws.Range("B3").Select()
